{"js": "// \"add a techniques sheet too\"\n// In the Kicho (basic forms) techniques table on the Orange Belt sheet,\n// fix up two of the technique names:\n//   \"Kicho E Jang\"   -> \"Kicho I Jang\"\n//   \"Kicho Sam Jan\"  -> \"Kicho Sam Jang\"\n// Note: the table cells use a non-breaking space (U+00A0) between words,\n// not a regular space, so we reproduce that exactly in both the search\n// and replacement text.\n\nconst nbsp = \"\\u00A0\";\nconst body = context.document.body;\n\n// 1) \"Kicho E Jang\" -> \"Kicho I Jang\"\nconst findKichoE = `Kicho${nbsp}E${nbsp}Jang`;\nconst replaceKichoI = `Kicho${nbsp}I${nbsp}Jang`;\nconst results1 = body.search(findKichoE, { matchCase: true, matchWholeWord: true });\nresults1.load(\"items\");\nawait context.sync();\n\nif (results1.items.length > 0) {\n  results1.items[0].insertText(replaceKichoI, \"Replace\");\n}\n\n// 2) \"Kicho Sam Jan\" -> \"Kicho Sam Jang\"\nconst findKichoSamJan = `Kicho${nbsp}Sam${nbsp}Jan`;\nconst replaceKichoSamJang = `Kicho${nbsp}Sam${nbsp}Jang`;\nconst results2 = body.search(findKichoSamJan, { matchCase: true, matchWholeWord: true });\nresults2.load(\"items\");\nawait context.sync();\n\nif (results2.items.length > 0) {\n  results2.items[0].insertText(replaceKichoSamJang, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# \"add a techniques sheet too\"\n# In the Kicho (basic forms) techniques table on the Orange Belt sheet,\n# fix up two of the technique names:\n#   \"Kicho E Jang\"   -> \"Kicho I Jang\"\n#   \"Kicho Sam Jan\"  -> \"Kicho Sam Jang\"\n# Note: the table cells use a non-breaking space (U+00A0) between words,\n# not a regular space, so we reproduce that exactly in both the search\n# and replacement text.\n\n$d = $word.ActiveDocument\n$nbsp = [char]0x00A0\n\n# 1) \"Kicho E Jang\" -> \"Kicho I Jang\"\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Replacement.ClearFormatting()\n$rng1.Find.Execute(\n    \"Kicho${nbsp}E${nbsp}Jang\",  # FindText\n    $false,                      # MatchCase\n    $true,                       # MatchWholeWord\n    $false,                      # MatchWildcards\n    $false,                      # MatchSoundsLike\n    $false,                      # MatchAllWordForms\n    $true,                       # Forward\n    1,                           # Wrap (wdFindContinue)\n    $false,                      # Format\n    \"Kicho${nbsp}I${nbsp}Jang\",  # ReplaceWith\n    2                            # Replace (wdReplaceAll)\n)\n\n# 2) \"Kicho Sam Jan\" -> \"Kicho Sam Jang\"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Replacement.ClearFormatting()\n$rng2.Find.Execute(\n    \"Kicho${nbsp}Sam${nbsp}Jan\",   # FindText\n    $false,                        # MatchCase\n    $true,                         # MatchWholeWord (avoid matching inside \"Jang\")\n    $false,                        # MatchWildcards\n    $false,                        # MatchSoundsLike\n    $false,                        # MatchAllWordForms\n    $true,                         # Forward\n    1,                             # Wrap (wdFindContinue)\n    $false,                        # Format\n    \"Kicho${nbsp}Sam${nbsp}Jang\",  # ReplaceWith\n    2                              # Replace (wdReplaceAll)\n)\n"}
